$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"
$mdDisplay  = "a1f78878-6f34-4aba-8c6b-09ecfc950b78.md"
$mdTarget   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/93dcd821e1e1ba944c3da0b0b83948cf05d958de/e2e/a1f78878-6f34-4aba-8c6b-09ecfc950b78.md"

# Column width (in "character" units) that rounds, once written back out as OOXML
# column width (points-ish units), to the widened values used by the edited file.
$wideColWidth   = 29.1666666666667   # -> stored width ~30 (was ~17.22)
$veryWideWidth  = 39.1666666666667   # -> stored width 40 (was ~18.65 / 21.71)

# ---------------------------------------------------------------------------
# Overview sheet: status column text for zh-cn / de-de changes
# ---------------------------------------------------------------------------
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = $wideColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideColWidth

# ---------------------------------------------------------------------------
# zh-cn sheet: handback is now complete
# ---------------------------------------------------------------------------
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Columns.Item(3).ColumnWidth = $wideColWidth

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdTarget, [Type]::Missing, [Type]::Missing, $mdDisplay)
$wsZhCn.Range("I2").Font.Underline = $true
$wsZhCn.Range("I2").Font.Color = 15570276

$wsZhCn.Range("J2").Value = "a1f78878-6f34-4aba-8c6b-09ecfc950b78.06b6b84ab7e855dab6f8e6d421c29ab64621fba8.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-27 06:55:44"

$wsZhCn.Columns.Item(9).ColumnWidth  = $veryWideWidth
$wsZhCn.Columns.Item(10).ColumnWidth = $veryWideWidth

# ---------------------------------------------------------------------------
# de-de sheet: handback is now complete
# ---------------------------------------------------------------------------
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Columns.Item(3).ColumnWidth = $wideColWidth

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdTarget, [Type]::Missing, [Type]::Missing, $mdDisplay)
$wsDeDe.Range("I2").Font.Underline = $true
$wsDeDe.Range("I2").Font.Color = 15570276

$wsDeDe.Range("J2").Value = "a1f78878-6f34-4aba-8c6b-09ecfc950b78.06b6b84ab7e855dab6f8e6d421c29ab64621fba8.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-27 06:55:50"

$wsDeDe.Columns.Item(9).ColumnWidth  = $veryWideWidth
$wsDeDe.Columns.Item(10).ColumnWidth = $veryWideWidth
